$d = $word.ActiveDocument

# 1. Trim the trailing English review fragment from the technical
#    documentation review title (keep the trailing space after "940330").
#    The search/replace strings deliberately avoid literal double-quote
#    characters so Word's "smart quotes" autocorrect does not touch the
#    quote marks that remain untouched elsewhere in the run.
$d.Content.Find.Execute(
    "РС 940330 / Review of technical documentation """" on mv VOLGA",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "РС 940330 ",
    2)

# 2. Add ", Россия" before the postal code for the Saint-Petersburg address.
$d.Content.Find.Execute(
    "Дворцовая набережная, 8, Санкт-Петербург  191186",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Дворцовая набережная, 8, Санкт-Петербург, Россия  191186",
    2)

# 3. Update both occurrences of the Murmansk address (building "19" -> "19\1",
#    add ", Россия" and a trailing space).
$d.Content.Find.Execute(
    "ул. Карла Маркса, д. 19, Мурманск  193025",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "ул. Карла Маркса, д. 19\1, Мурманск, Россия  193025 ",
    2)

# 4. Add ", Россия" before the postal code for the Kaliningrad address.
$d.Content.Find.Execute(
    "Молочинского, д. 4, Калининград  236023",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Молочинского, д. 4, Калининград, Россия  236023",
    2)
